# Update expiry dates (column B) for the symbols sheet.
# Each changed expiry is shifted forward by 7 days (one week), reflecting
# the weekly option-expiry roll described in the commit message.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value  = 45498   # NIFTY      45491 -> 45498
$ws.Range("B3").Value  = 45505   # NIFTY      45498 -> 45505
$ws.Range("B8").Value  = 45503   # FINNIFTY   45496 -> 45503
$ws.Range("B9").Value  = 45510   # FINNIFTY   45503 -> 45510
$ws.Range("B10").Value = 45502   # MIDCPNIFTY 45495 -> 45502
$ws.Range("B11").Value = 45509   # MIDCPNIFTY 45502 -> 45509
